$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value could be misread as a number by Excel;
# force them to Text format before writing, then clear the formatting
# residue afterwards so no stray style index is left on the cell.
$forceTextCells = @("D5", "D6", "D7", "D11", "D16", "D19", "D21", "D22", "D27", "D30", "D31", "D33", "D36", "D38", "D42", "D44", "D45", "D46", "D47", "D49", "D50", "D51")
foreach ($c in $forceTextCells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range("D2").Value = '66.224.31'
$ws.Range("E2").Value = '  -0.80%  '
$ws.Range("D3").Value = '3.320.23'
$ws.Range("E3").Value = '  -1.09%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '585.96'
$ws.Range("E5").Value = '  +2.08%  '
$ws.Range("D6").Value = '181.26'
$ws.Range("E6").Value = '  -0.38%  '
$ws.Range("D7").Value = '0.651'
$ws.Range("E7").Value = '  +3.65%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("D9").Value = '3.320.21'
$ws.Range("E9").Value = '  -1.03%  '
$ws.Range("E10").Value = '  -2.32%  '
$ws.Range("D11").Value = '6.83'
$ws.Range("E11").Value = '  +2.46%  '
$ws.Range("E12").Value = '  -0.56%  '
$ws.Range("D13").Value = '3.898.22'
$ws.Range("E13").Value = '  -1.01%  '
$ws.Range("E14").Value = '  -2.43%  '
$ws.Range("D15").Value = '66.245.17'
$ws.Range("E15").Value = '  -0.89%  '
$ws.Range("D16").Value = '26.35'
$ws.Range("E16").Value = '  -2.28%  '
$ws.Range("E17").Value = '  -1.32%  '
$ws.Range("D18").Value = '3.314.82'
$ws.Range("E18").Value = '  -1.54%  '
$ws.Range("D19").Value = '425.26'
$ws.Range("E19").Value = '  -3.70%  '
$ws.Range("E20").Value = '  -2.65%  '
$ws.Range("D21").Value = '13.14'
$ws.Range("E21").Value = '  -3.44%  '
$ws.Range("D22").Value = '7.40'
$ws.Range("E22").Value = '  -2.71%  '
$ws.Range("E23").Value = '  -2.41%  '
$ws.Range("E24").Value = '  +0.25%  '
$ws.Range("E25").Value = '  +0.14%  '
$ws.Range("D26").Value = '3.463.87'
$ws.Range("E26").Value = '  -1.04%  '
$ws.Range("D27").Value = '0.515'
$ws.Range("E27").Value = '  -1.02%  '
$ws.Range("E28").Value = '  +4.89%  '
$ws.Range("E29").Value = '  -2.07%  '
$ws.Range("D30").Value = '8.97'
$ws.Range("E30").Value = '  -1.75%  '
$ws.Range("D31").Value = '0.995'
$ws.Range("E31").Value = '  -0.48%  '
$ws.Range("E32").Value = '  -1.71%  '
$ws.Range("D33").Value = '22.37'
$ws.Range("E33").Value = '  -2.26%  '
$ws.Range("E34").Value = '  +0.03%  '
$ws.Range("E35").Value = '  -2.54%  '
$ws.Range("D36").Value = '6.60'
$ws.Range("E36").Value = '  -3.29%  '
$ws.Range("E37").Value = '  -3.61%  '
$ws.Range("D38").Value = '160.35'
$ws.Range("E38").Value = '  -0.80%  '
$ws.Range("E39").Value = '  -3.36%  '
$ws.Range("D40").Value = '2.870.67'
$ws.Range("E40").Value = '  +1.47%  '
$ws.Range("E41").Value = '  +0.27%  '
$ws.Range("D42").Value = '26.36'
$ws.Range("E42").Value = '  -6.04%  '
$ws.Range("E43").Value = '  -2.79%  '
$ws.Range("D44").Value = '0.758'
$ws.Range("E44").Value = '  -5.27%  '
$ws.Range("D45").Value = '39.79'
$ws.Range("E45").Value = '  -1.66%  '
$ws.Range("D46").Value = '0.0661'
$ws.Range("E46").Value = '  -1.15%  '
$ws.Range("D47").Value = '5.91'
$ws.Range("E47").Value = '  -4.96%  '
$ws.Range("E48").Value = '  -2.05%  '
$ws.Range("D49").Value = '23.13'
$ws.Range("E49").Value = '  -4.93%  '
$ws.Range("D50").Value = '312.82'
$ws.Range("E50").Value = '  -4.51%  '
$ws.Range("D51").Value = '0.0272'
$ws.Range("E51").Value = '  -0.39%  '

foreach ($c in $forceTextCells) {
    $ws.Range($c).ClearFormats()
}
